$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.811.21'
$ws.Range("E2").Value = '  -4.05%  '
$ws.Range("D3").Value = '3.133.61'
$ws.Range("E3").Value = '  -3.60%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.39'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.90%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '3.133.05'
$ws.Range("E8").Value = '  -3.65%  '
$ws.Range("E9").Value = '  -4.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.150'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.35'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -8.39%  '
$ws.Range("E12").Value = '  -5.97%  '
$ws.Range("E13").Value = '  -8.22%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.27'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -9.65%  '
$ws.Range("D15").Value = '3.644.82'
$ws.Range("E15").Value = '  -3.79%  '
$ws.Range("E16").Value = '  +0.55%  '
$ws.Range("D17").Value = '63.881.45'
$ws.Range("E17").Value = '  -4.10%  '
$ws.Range("D18").Value = '3.131.76'
$ws.Range("E18").Value = '  -3.76%  '
$ws.Range("E19").Value = '  -8.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '474.94'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -6.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.70'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.56%  '
$ws.Range("E22").Value = '  -6.68%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.76'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.53'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -8.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.84'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.88%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("E27").Value = '  -4.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.39'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -8.24%  '
$ws.Range("E29").Value = '  -9.45%  '
$ws.Range("E30").Value = '  -3.55%  '
$ws.Range("E31").Value = '  -12.46%  '
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("E33").Value = '  -6.42%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.07'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.94'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.44'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.78%  '
$ws.Range("D38").Value = '0.0₃0727'
$ws.Range("E38").Value = '  -7.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '458.27'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.20%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.92'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -13.56%  '
$ws.Range("E41").Value = '  -7.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.39'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.95%  '
$ws.Range("E43").Value = '  -7.79%  '
$ws.Range("D44").Value = '2.832.98'
$ws.Range("E44").Value = '  -5.10%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.264'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -9.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -10.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.31'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -9.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '118.94'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.86%  '
